$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "80-16="
$t.Cell(1,2).Range.Text = "63-49="
$t.Cell(1,3).Range.Text = "28+43="
$t.Cell(1,4).Range.Text = "47+47="
$t.Cell(1,5).Range.Text = "83-17="
$t.Cell(2,1).Range.Text = "59+12="
$t.Cell(2,2).Range.Text = "35+8="
$t.Cell(2,3).Range.Text = "3+89="
$t.Cell(2,4).Range.Text = "25+19="
$t.Cell(2,5).Range.Text = "3+38="
$t.Cell(3,1).Range.Text = "53-28="
$t.Cell(3,2).Range.Text = "92-23="
$t.Cell(3,3).Range.Text = "76+5="
$t.Cell(3,4).Range.Text = "81-2="
$t.Cell(3,5).Range.Text = "45-37="
$t.Cell(4,1).Range.Text = "8+18="
$t.Cell(4,2).Range.Text = "63-9="
$t.Cell(4,3).Range.Text = "68+6="
$t.Cell(4,4).Range.Text = "52-9="
$t.Cell(4,5).Range.Text = "94-46="
$t.Cell(5,1).Range.Text = "18+73="
$t.Cell(5,2).Range.Text = "25-18="
$t.Cell(5,3).Range.Text = "79+12="
$t.Cell(5,4).Range.Text = "94-26="
$t.Cell(5,5).Range.Text = "60-16="
$t.Cell(6,1).Range.Text = "49+6="
$t.Cell(6,2).Range.Text = "97-19="
$t.Cell(6,3).Range.Text = "85+7="
$t.Cell(6,4).Range.Text = "9+29="
$t.Cell(6,5).Range.Text = "49+8="
$t.Cell(7,1).Range.Text = "89+4="
$t.Cell(7,2).Range.Text = "15+17="
$t.Cell(7,3).Range.Text = "49+34="
$t.Cell(7,4).Range.Text = "18+9="
$t.Cell(7,5).Range.Text = "42+19="
$t.Cell(8,1).Range.Text = "30-24="
$t.Cell(8,2).Range.Text = "8+54="
$t.Cell(8,3).Range.Text = "44+19="
$t.Cell(8,4).Range.Text = "77+16="
$t.Cell(8,5).Range.Text = "48+5="
$t.Cell(9,1).Range.Text = "68-39="
$t.Cell(9,2).Range.Text = "17+14="
$t.Cell(9,3).Range.Text = "51-22="
$t.Cell(9,4).Range.Text = "96-9="
$t.Cell(9,5).Range.Text = "71-4="
$t.Cell(10,1).Range.Text = "29+67="
$t.Cell(10,2).Range.Text = "15+78="
$t.Cell(10,3).Range.Text = "8+4="
$t.Cell(10,4).Range.Text = "75-69="
$t.Cell(10,5).Range.Text = "22-14="
$t.Cell(11,1).Range.Text = "34-6="
$t.Cell(11,2).Range.Text = "64-39="
$t.Cell(11,3).Range.Text = "32+9="
$t.Cell(11,4).Range.Text = "15+46="
$t.Cell(11,5).Range.Text = "90-62="
$t.Cell(12,1).Range.Text = "58+9="
$t.Cell(12,2).Range.Text = "62-45="
$t.Cell(12,3).Range.Text = "14+8="
$t.Cell(12,4).Range.Text = "58+39="
$t.Cell(12,5).Range.Text = "52-38="
$t.Cell(13,1).Range.Text = "53-34="
$t.Cell(13,2).Range.Text = "35-18="
$t.Cell(13,3).Range.Text = "84-46="
$t.Cell(13,4).Range.Text = "73-6="
$t.Cell(13,5).Range.Text = "68-19="
$t.Cell(14,1).Range.Text = "70-51="
$t.Cell(14,2).Range.Text = "24+39="
$t.Cell(14,3).Range.Text = "60-53="
$t.Cell(14,4).Range.Text = "26+69="
$t.Cell(14,5).Range.Text = "49+25="
$t.Cell(15,1).Range.Text = "90-44="
$t.Cell(15,2).Range.Text = "92-23="
$t.Cell(15,3).Range.Text = "9+75="
$t.Cell(15,4).Range.Text = "9+74="
$t.Cell(15,5).Range.Text = "35-26="
$t.Cell(16,1).Range.Text = "90-35="
$t.Cell(16,2).Range.Text = "84-59="
$t.Cell(16,3).Range.Text = "94-16="
$t.Cell(16,4).Range.Text = "47-9="
$t.Cell(16,5).Range.Text = "28+68="
$t.Cell(17,1).Range.Text = "78+9="
$t.Cell(17,2).Range.Text = "76-29="
$t.Cell(17,3).Range.Text = "51-18="
$t.Cell(17,4).Range.Text = "27+46="
$t.Cell(17,5).Range.Text = "38+37="
$t.Cell(18,1).Range.Text = "57+36="
$t.Cell(18,2).Range.Text = "2+89="
$t.Cell(18,3).Range.Text = "57+7="
$t.Cell(18,4).Range.Text = "93-39="
$t.Cell(18,5).Range.Text = "60-5="
$t.Cell(19,1).Range.Text = "66+25="
$t.Cell(19,2).Range.Text = "43-4="
$t.Cell(19,3).Range.Text = "26+59="
$t.Cell(19,4).Range.Text = "25+17="
$t.Cell(19,5).Range.Text = "90-41="
$t.Cell(20,1).Range.Text = "47-18="
$t.Cell(20,2).Range.Text = "53-49="
$t.Cell(20,3).Range.Text = "37+57="
$t.Cell(20,4).Range.Text = "58+24="
$t.Cell(20,5).Range.Text = "26+69="
